# Update "想去人数" (F) counts and clear "已售罄" (sold out) status on G for
# rows where tickets are available again, matching the refreshed data pull.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F3").Value = 1298
$ws.Range("F4").Value = 13010
$ws.Range("F5").Value = 742
$ws.Range("F7").Value = 319
$ws.Range("F10").Value = 1889
$ws.Range("F13").Value = 5990
$ws.Range("F15").Value = 216
$ws.Range("F17").Value = 355
$ws.Range("F20").Value = 143
$ws.Range("F22").Value = 30
$ws.Range("F23").Value = 227
$ws.Range("F25").Value = 1318
$ws.Range("F26").Value = 349

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F6").Value = 168
$ws.Range("F11").Value = 368

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 888
$ws.Range("F3").Value = 4306
$ws.Range("G3").Value = 0

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 888
$ws.Range("F6").Value = 1298
$ws.Range("F7").Value = 13010
$ws.Range("F9").Value = 742
$ws.Range("F10").Value = 4306
$ws.Range("G10").Value = 0
$ws.Range("F12").Value = 319
$ws.Range("F15").Value = 1889
$ws.Range("F18").Value = 5995
$ws.Range("G20").Value = 0
$ws.Range("F21").Value = 216
$ws.Range("F22").Value = 168
$ws.Range("F23").Value = 168
$ws.Range("F29").Value = 368
$ws.Range("F30").Value = 355
$ws.Range("F34").Value = 143
$ws.Range("F36").Value = 30
$ws.Range("F38").Value = 227
$ws.Range("F42").Value = 1318
$ws.Range("F44").Value = 349
